# Add a new "libraryProtocol" column (J) to the library sheet and fill it
# in for every existing data row with the protocol code "E7420L".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new column.
$ws.Range("J1").Value = "libraryProtocol"
$ws.Range("J1").Font.Name = "Arial"
$ws.Range("J1").Font.Size = 10

# Fill the protocol value down for all existing data rows (2-19).
$ws.Range("J2:J19").Value = "E7420L"
$ws.Range("J2:J19").Font.Name = "Arial"
$ws.Range("J2:J19").Font.Size = 11

# Row 1's height settles to match the rest of the sheet once the new
# column is populated.
$ws.Rows.Item(1).RowHeight = 13.8

# A handful of extra (still empty) rows below the table picked up the same
# row height while the selection was being dragged down.
for ($r = 20; $r -le 27; $r++) {
    $ws.Rows.Item($r).RowHeight = 13.8
}

# Leave the selection where the user ended up after filling the column.
$ws.Range("J20:J27").Select()
